$wb = $excel.ActiveWorkbook

# Both "VENTAS POR GRUPO" (cols C:R) and "VENTA MENSUAL" (cols C:G) sheets get
# a new client row inserted right before "PUCO TOAPANTA MARCO ANTONIO" (row 43),
# for new client "PACHAR TAPIA ELIANA DE LOS ANGELES" with zeroed metrics. This
# pushes every subsequent row (and the trailing totals row) down by one.

$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(43).Insert()
$ws1.Range("A43").Value = "LINDAO ZUÑIGA BRYAN JOSE"
$ws1.Range("B43").Value = "PACHAR TAPIA ELIANA DE LOS ANGELES"
$ws1.Range("C43:R43").Value = 0

# Totals row (now row 60) text counters: "x de 57" -> "x de 58"
$ws1.Range("C60:L60").Value = "0 de 58"
$ws1.Range("M60").Value = "1 de 58"
$ws1.Range("N60:R60").Value = "0 de 58"

$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(43).Insert()
$ws2.Range("A43").Value = "LINDAO ZUÑIGA BRYAN JOSE"
$ws2.Range("B43").Value = "PACHAR TAPIA ELIANA DE LOS ANGELES"
$ws2.Range("C43:G43").Value = 0
